# Update the daily Scores sheet with the latest computed values.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 74 (2025-02-19, abs_activity)
$ws.Range("D74").Value = 10
$ws.Range("F74").Value = 20

# Row 75 (2025-02-19, rel_activity)
$ws.Range("D75").Value = 5.536250536250536
$ws.Range("F75").Value = 5.536250536250536

# Row 78 (2025-02-20, abs_activity)
$ws.Range("C78").Value = 9.543199674630392
$ws.Range("D78").Value = 6.086067018668716
$ws.Range("F78").Value = 15.62926669329911

# Row 80 (2025-02-20, abs_sleep)
$ws.Range("C80").Value = 5.266666666666666
$ws.Range("D80").Value = 10
$ws.Range("F80").Value = 15.26666666666667

# Row 81 (2025-02-20, rel_sleep)
$ws.Range("D81").Value = 9.287132265008371
$ws.Range("F81").Value = 9.287132265008371
